$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple price updates (column D only)
$ws.Range("D2").Value = "242.69"
$ws.Range("D3").Value = "23.03"
$ws.Range("D4").Value = "5.412"
$ws.Range("D5").Value = "0.05899"
$ws.Range("D7").Value = "6.530"
$ws.Range("D8").Value = "0.8086"
$ws.Range("D9").Value = "0.9393"
$ws.Range("D11").Value = "0.07403"
$ws.Range("D12").Value = "0.03288"
$ws.Range("D13").Value = "0.03066"
$ws.Range("D14").Value = "0.09347"
$ws.Range("D15").Value = "3.848"
$ws.Range("D16").Value = "0.001572"
$ws.Range("D17").Value = "0.04690"

# Rows 18-24: data shifted down by one row, with new "One" row inserted at 18
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "0.0005921"
$ws.Range("E18").Value = "17OneONE"

$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "0.006019"
$ws.Range("E19").Value = "18TigerCashTCH"

$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "0.001258"
$ws.Range("E20").Value = "19BitKanKAN"

$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").Value = "0.004905"
$ws.Range("E21").Value = "20HotbitTokenHTB"

$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "0.00006802"
$ws.Range("E22").Value = "21NitroExNTX"

$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "3.562"
$ws.Range("E23").Value = "22LEOLEO"

$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "2.124"
$ws.Range("E24").Value = "23BTSETokenBTSE"

# Rows 40-49 price updates
$ws.Range("D40").Value = "0.03961"
$ws.Range("D41").Value = "0.006178"
$ws.Range("D42").Value = "0.1072"
$ws.Range("D43").Value = "0.002571"
$ws.Range("D44").Value = "0.009497"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"
$ws.Range("D47").Value = "0.6701"
$ws.Range("D48").Value = "0.002329"
$ws.Range("D49").Value = "0.00002100"
